$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# New header cell M1 = "Event", matching the formatting of the other
# header cells in row 1 (bold, bordered, centered - same style as L1).
$ws.Cells.Item(1, 12).Copy() | Out-Null
$ws.Cells.Item(1, 13).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 13).Value = "Event"

# New (blank) data cells M2:M12 under the new "Event" column.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Style = "Normal"
}
